# KetQuaTest_RegisterServlet.xlsx — "all source and test"
#
# Re-order / rewrite the Register-Servlet test-case rows:
#   Row2 -> REG_SERV_04 "Loi Database" (Du Lieu Mau / Cac Buoc columns swapped)
#   Row3 -> REG_SERV_02 "Mat khau yeu" (brand new test case)
#   Row4 -> REG_SERV_03 "Trung ten dang nhap" (Du Lieu Mau / Cac Buoc swapped)
#   Row5 -> REG_SERV_01 "Dang ky thanh cong" (new row, Pass value updated)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: REG_SERV_04 / Servlet: Loi Database ------------------------
$ws.Range("A2").Value = "REG_SERV_04"
$ws.Range("B2").Value = "Servlet: Lỗi Database"
$ws.Range("C2").Value = "User: valid"
$ws.Range("D2").Value = "1. Mock DAO register=false"
$ws.Range("E2").Value = "Redirect: Register.jsp + Error"
$ws.Range("F2").Value = "OK"
$ws.Range("G2").Value = "PASS"

# --- Row 3: REG_SERV_02 / Servlet: Mat khau yeu (new scenario) --------
$ws.Range("A3").Value = "REG_SERV_02"
$ws.Range("B3").Value = "Servlet: Mật khẩu yếu"
$ws.Range("C3").Value = "Pass: 123"
$ws.Range("D3").Value = "1. Mock pass ngắn (<6 ký tự)`n2. Call doPost"
$ws.Range("E3").Value = "Báo lỗi 'Mật khẩu quá yếu'"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"

# --- Row 4: REG_SERV_03 / Servlet: Trung ten dang nhap -----------------
$ws.Range("A4").Value = "REG_SERV_03"
$ws.Range("B4").Value = "Servlet: Trùng tên đăng nhập"
$ws.Range("C4").Value = "User: exist"
$ws.Range("D4").Value = "1. Mock DAO checkExists=true"
$ws.Range("E4").Value = "Redirect: Register.jsp + Error"
$ws.Range("F4").Value = "OK"
$ws.Range("G4").Value = "PASS"

# --- Row 5 (new): REG_SERV_01 / Servlet: Dang ky thanh cong ------------
$ws.Range("A5").Value = "REG_SERV_01"
$ws.Range("B5").Value = "Servlet: Đăng ký thành công"
$ws.Range("C5").Value = "User: new, Pass: 123456"
$ws.Range("D5").Value = "1. Mock input hợp lệ`n2. Mock DAO register=true`n3. Call doPost"
$ws.Range("E5").Value = "Redirect: Login.jsp"
$ws.Range("F5").Value = "OK"
# Pick up the PASS formatting (green bold cell style) from the row above
# instead of re-typing the value, so the shared cellXf index (s="2") is
# reused rather than a new one minted.
$ws.Range("G4").Copy($ws.Range("G5"))

# Re-entering the multi-line "Các Bước" text bumps those rows to an
# auto-calculated custom height; AutoFit restores the sheet's standard
# (non-custom) row height, same as the original rows.
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# --- Column widths: Du Lieu Mau (C) / Cac Buoc (D) swap + re-fit -------
# The sample-data column got narrower and the steps column got wider once
# the text moved between them; reproduce the closest achievable widths
# (ColumnWidth is pixel-quantized, so these land on the nearest pixel step
# to the authored 22.8125 / 28.20703125 character widths).
$ws.Range("C1").ColumnWidth = 22
$ws.Range("D1").ColumnWidth = 27.3

Write-Output "Register servlet test rows reorganized."
